$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column P (2021 data) holds the same formatting/placeholder pattern that the
# new column Q (2022 data) needs, so copy P3:P25 -> Q3:Q25 to bring across
# styles (borders, number formats, "..." placeholders) in one shot.
$src = $ws.Range("P3:P25")
$dst = $ws.Range("Q3:Q25")
$src.Copy($dst)

# Now overwrite the copied values with the actual 2022 figures.
$ws.Range("Q4").Value = 2022
$ws.Range("Q5").Value = 8725
$ws.Range("Q7").Value = 8347
$ws.Range("Q8").Value = 378

# Update the selected cell to Q3, matching the new sheet view selection.
$ws.Range("Q3").Select() | Out-Null
